$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# This quarterly update ("Actualización desde MV -datos-") revises the
# interest-related figures (columns J, M, Q, W, X) for the two most
# recent existing rows (2021 Q1 and Q2) and appends a brand-new row for
# 2021 Q3 ("01-07-2021").

# --- Revise row 66 (01-01-2021): 6980 -> 7013 ---
$ws.Cells.Item(66, 10).Value = 7013   # J
$ws.Cells.Item(66, 13).Value = 7013   # M
$ws.Cells.Item(66, 17).Value = -7013  # Q
$ws.Cells.Item(66, 23).Value = 7013   # W
$ws.Cells.Item(66, 24).Value = -7013  # X

# --- Revise row 67 (01-04-2021): 6411 -> 6444 ---
$ws.Cells.Item(67, 10).Value = 6444   # J
$ws.Cells.Item(67, 13).Value = 6444   # M
$ws.Cells.Item(67, 17).Value = -6444  # Q
$ws.Cells.Item(67, 23).Value = 6444   # W
$ws.Cells.Item(67, 24).Value = -6444  # X

# --- Append new row 68 (01-07-2021) ---
$row = 68

# Force column A to be stored as text so the date-like label
# "01-07-2021" isn't auto-converted to a date serial number, then drop
# the temporary number format so the cell keeps the default style.
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "01-07-2021"
$ws.Cells.Item($row, 1).ClearFormats()

for ($col = 2; $col -le 24; $col++) {
    $ws.Cells.Item($row, $col).Value = 0
}

$ws.Cells.Item($row, 10).Value = 5866   # J - Gastos
$ws.Cells.Item($row, 13).Value = 5866   # M - Intereses
$ws.Cells.Item($row, 17).Value = -5866  # Q - Resultado operativo bruto
$ws.Cells.Item($row, 23).Value = 5866   # W - Total gastos
$ws.Cells.Item($row, 24).Value = -5866  # X - Préstamo o endeudamiento neto
